$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("A3").Value = "a"
$ws.Range("B3").Value = 0.0
